$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate stats now that trade #35 has closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.74
$summary.Range("B4").Value = 0.73
$summary.Range("B5").Value = 0.42
$summary.Range("B6").Value = 35
$summary.Range("B7").Value = 14
$summary.Range("B9").Value = 40

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) picks up the new trade
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.74
$status.Range("D4").Value = 35
$status.Range("E4").Value = 0.73
$status.Range("F4").Value = 0.74
$status.Range("G4").Value = 40

# ---------------------------------------------------------------------------
# Helper: append the newly-closed trade #35 to a trade-log sheet (row 36)
# Text columns are forced to the "Text" number format before the write so
# Excel stores them as literal strings (e.g. dates/times) instead of
# auto-converting to date/time serials, then the style is put back to
# Normal so no stray formatting is left behind on the cell.
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

function Add-Trade35Row($ws) {
    $ws.Range("A36").Value = 35
    Set-TextCell $ws "B36" "2026-02-17"
    Set-TextCell $ws "C36" "12:38:49"
    Set-TextCell $ws "D36" "MarketMaking"
    Set-TextCell $ws "E36" "UP"
    $ws.Range("F36").Value = 0.821782
    $ws.Range("G36").Value = 0.91
    Set-TextCell $ws "H36" "CLOSED"
    $ws.Range("I36").Value = 10.7349
    $ws.Range("J36").Value = 0.09
    $ws.Range("K36").Value = 100.74
    $ws.Range("L36").Value = 0
    $ws.Range("M36").Value = 0
    $ws.Range("N36").Value = 0.6
    Set-TextCell $ws "O36" "Normal spread capture: 19600 bps"
    Set-TextCell $ws "P36" "early_exit"
    $ws.Range("Q36").Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet: append trade #35
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade35Row $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet: append the same trade #35
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade35Row $marketMaking
